# Auto-generated edit script applying the cell-value diff to cryptos worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.613.85"
$ws.Range("E2").Value = "  +4.07%  "
$ws.Range("D3").Value = "2.267.33"
$ws.Range("E3").Value = "  +1.36%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'230.10"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = "  -0.08%  "
$ws.Range("D7").Value = "'63.23"
$ws.Range("E7").Value = "  +4.45%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.421"
$ws.Range("E9").Value = "  +3.82%  "
$ws.Range("D10").Value = "'0.0988"
$ws.Range("E10").Value = "  +9.03%  "
$ws.Range("D11").Value = "'57.39"
$ws.Range("E11").Value = "  -1.07%  "
$ws.Range("D12").Value = "'25.62"
$ws.Range("E12").Value = "  +12.32%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "2.604.89"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "'15.56"
$ws.Range("E15").Value = "  -0.10%  "
$ws.Range("D16").Value = "'5.85"
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "'0.806"
$ws.Range("E17").Value = "  +0.14%  "
$ws.Range("D18").Value = "2.288.39"
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("D19").Value = "43.550.34"
$ws.Range("E19").Value = "  +4.05%  "
$ws.Range("D20").Value = "0.0₃0970"
$ws.Range("E20").Value = "  +5.93%  "
$ws.Range("D21").Value = "'72.73"
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("D22").Value = "'6.05"
$ws.Range("E22").Value = "  -1.84%  "
$ws.Range("D23").Value = "'247.25"
$ws.Range("E23").Value = "  -0.66%  "
$ws.Range("E24").Value = "  +0.23%  "
$ws.Range("D25").Value = "'2.47"
$ws.Range("E25").Value = "  +3.72%  "
$ws.Range("E26").Value = "  +0.82%  "
$ws.Range("D27").Value = "'9.85"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").Value = "'171.28"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "'0.136"
$ws.Range("E29").Value = "  -3.64%  "
$ws.Range("D30").Value = "'20.43"
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").Value = "'1.44"
$ws.Range("E31").Value = "  +1.41%  "
$ws.Range("E32").Value = "  +8.93%  "
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "'0.0684"
$ws.Range("E34").Value = "  +3.81%  "
$ws.Range("D35").Value = "'5.05"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "'4.66"
$ws.Range("E36").Value = "  -1.38%  "
$ws.Range("D37").Value = "'3.82"
$ws.Range("E37").Value = "  +5.05%  "
$ws.Range("D38").Value = "'6.64"
$ws.Range("E38").Value = "  +0.20%  "
$ws.Range("E39").Value = "  -4.28%  "
$ws.Range("D40").Value = "'0.0245"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").Value = "'8.32"
$ws.Range("E42").Value = "  -4.31%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D43").Value = "'10.39"
$ws.Range("E43").Value = "  +17.90%  "
$ws.Range("B44").Value = "Cronos"
$ws.Range("C44").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D44").Value = "'0.0958"
$ws.Range("E44").Value = "  -0.78%  "
$ws.Range("B45").Value = "TrustWalletToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D45").Value = "'1.20"
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("B46").Value = "InjectiveProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D46").Value = "'16.93"
$ws.Range("E46").Value = "  +0.96%  "
$ws.Range("B47").Value = "TerraClassic"
$ws.Range("C47").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D47").Value = "'0.000211"
$ws.Range("E47").Value = "  -15.17%  "
$ws.Range("D48").Value = "'96.56"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("D49").Value = "1.468.35"
$ws.Range("E49").Value = "  -0.71%  "
$ws.Range("B50").Value = "FTXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D50").Value = "'4.35"
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("E51").Value = "  +0.67%  "
